$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Save" header column (H1), matching the formatting of the other
# header cells (e.g. G1: bold, centered, bordered) by copying its format.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# New data value for the Save column (H2)
$ws.Range("H2").Value = 1
